$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.030141
$ws.Range("H2").Value = 0.090423
$ws.Range("I2").Value = 0.2387784203438168
$ws.Range("J2").Value = 0.2387784203438169
$ws.Range("M2").Value = 0.06624833333333334
$ws.Range("N2").Value = 0.198745
$ws.Range("Q2").Value = 0.001996791015
$ws.Range("R2").Value = 0.017971119135
$ws.Range("S2").Value = 0.2387784203438168
$ws.Range("T2").Value = 0.2387784203438169

# Row 3
$ws.Range("H3").Value = 0.08290500000000001
$ws.Range("I3").Value = 0.2189257704190763
$ws.Range("J3").Value = 0.2189257704190763
$ws.Range("M3").Value = 0.06624833333333334
$ws.Range("N3").Value = 0.198745
$ws.Range("Q3").Value = 0.001830772691666667
$ws.Range("R3").Value = 0.016476954225
$ws.Range("S3").Value = 0.2189257704190763
$ws.Range("T3").Value = 0.2189257704190763

# Row 4
$ws.Range("G4").Value = 0.068454
$ws.Range("H4").Value = 0.205362
$ws.Range("I4").Value = 0.5422958092371069
$ws.Range("J4").Value = 0.5422958092371069
$ws.Range("M4").Value = 0.06624833333333334
$ws.Range("N4").Value = 0.198745
$ws.Range("Q4").Value = 0.00453496341
$ws.Range("R4").Value = 0.04081467069
$ws.Range("S4").Value = 0.5422958092371069
$ws.Range("T4").Value = 0.5422958092371069
